$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header team labels (K1<->N1, L1<->O1) ---
$k1 = $ws.Range("K1").Value2
$l1 = $ws.Range("L1").Value2
$n1 = $ws.Range("N1").Value2
$o1 = $ws.Range("O1").Value2
$ws.Range("K1").Value = $n1
$ws.Range("L1").Value = $o1
$ws.Range("N1").Value = $k1
$ws.Range("O1").Value = $l1

# --- Swap header colors so branding follows the team ---
$k1int = $ws.Range("K1").Interior.Color
$k1fnt = $ws.Range("K1").Font.Color
$l1int = $ws.Range("L1").Interior.Color
$l1fnt = $ws.Range("L1").Font.Color
$n1int = $ws.Range("N1").Interior.Color
$n1fnt = $ws.Range("N1").Font.Color
$o1int = $ws.Range("O1").Interior.Color
$o1fnt = $ws.Range("O1").Font.Color
$ws.Range("K1").Interior.Color = $n1int
$ws.Range("K1").Font.Color = $n1fnt
$ws.Range("L1").Interior.Color = $o1int
$ws.Range("L1").Font.Color = $o1fnt
$ws.Range("N1").Interior.Color = $k1int
$ws.Range("N1").Font.Color = $k1fnt
$ws.Range("O1").Interior.Color = $l1int
$ws.Range("O1").Font.Color = $l1fnt

# --- Swap location sub-labels (K2<->N2) ---
$k2 = $ws.Range("K2").Value2
$n2 = $ws.Range("N2").Value2
$ws.Range("K2").Value = $n2
$ws.Range("N2").Value = $k2

# --- Forecast numeric values (rows 3-26) ---
$ws.Range("K3").Value = 0.1324999320622214
$ws.Range("N3").Value = 0.4892509403539777
$ws.Range("B4").Value = 0.03486570345495707
$ws.Range("E4").Value = 0.9659600409376649
$ws.Range("H4").Value = 0.3494513174963016
$ws.Range("K4").Value = 0.7303307061449437
$ws.Range("N4").Value = 0.06397139186840474
$ws.Range("B5").Value = 1.767243353768706
$ws.Range("E5").Value = 47.05277035946656
$ws.Range("H5").Value = 23.84324356428177
$ws.Range("K5").Value = 9.676876894715923
$ws.Range("N5").Value = 3.129806362736982
$ws.Range("B6").Value = 0.9945311999999999
$ws.Range("C6").Value = 0.0036498
$ws.Range("E6").Value = 0.5700148
$ws.Range("F6").Value = 0.3672224
$ws.Range("H6").Value = 0.9145027999999999
$ws.Range("I6").Value = 0.0642486
$ws.Range("K6").Value = 0.776268
$ws.Range("L6").Value = 0.1992392
$ws.Range("N6").Value = 0.0075164
$ws.Range("O6").Value = 0.9899476
$ws.Range("B7").Value = 35.1190718
$ws.Range("C7").Value = 13.82852
$ws.Range("E7").Value = 28.0165104
$ws.Range("F7").Value = 25.3710078
$ws.Range("H7").Value = 28.1554178
$ws.Range("I7").Value = 17.1412582
$ws.Range("K7").Value = 24.444969
$ws.Range("L7").Value = 16.2451552
$ws.Range("N7").Value = 11.6851326
$ws.Range("O7").Value = 34.4049178
$ws.Range("K8").Value = 12
$ws.Range("L8").Value = 5
$ws.Range("N8").Value = 6
$ws.Range("O8").Value = 22
$ws.Range("K9").Value = 14
$ws.Range("L9").Value = 7
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 25
$ws.Range("K10").Value = 17
$ws.Range("L10").Value = 8
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 27
$ws.Range("K11").Value = 19
$ws.Range("L11").Value = 10
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 28
$ws.Range("K12").Value = 19
$ws.Range("L12").Value = 11
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 30
$ws.Range("K13").Value = 21
$ws.Range("L13").Value = 12
$ws.Range("N13").Value = 9
$ws.Range("O13").Value = 31
$ws.Range("K14").Value = 22
$ws.Range("L14").Value = 13
$ws.Range("N14").Value = 9
$ws.Range("O14").Value = 32
$ws.Range("K15").Value = 22
$ws.Range("L15").Value = 14
$ws.Range("N15").Value = 9
$ws.Range("O15").Value = 33
$ws.Range("K16").Value = 24
$ws.Range("L16").Value = 15
$ws.Range("N16").Value = 9
$ws.Range("O16").Value = 34
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = 16
$ws.Range("N17").Value = 10
$ws.Range("O17").Value = 35
$ws.Range("K18").Value = 26
$ws.Range("L18").Value = 17
$ws.Range("N18").Value = 12
$ws.Range("O18").Value = 36
$ws.Range("K19").Value = 26
$ws.Range("L19").Value = 18
$ws.Range("N19").Value = 13
$ws.Range("O19").Value = 36
$ws.Range("K20").Value = 28
$ws.Range("L20").Value = 19
$ws.Range("N20").Value = 13
$ws.Range("O20").Value = 37
$ws.Range("K21").Value = 28
$ws.Range("L21").Value = 20
$ws.Range("N21").Value = 14
$ws.Range("O21").Value = 38
$ws.Range("K22").Value = 29
$ws.Range("L22").Value = 21
$ws.Range("N22").Value = 16
$ws.Range("O22").Value = 39
$ws.Range("K23").Value = 31
$ws.Range("L23").Value = 23
$ws.Range("N23").Value = 16
$ws.Range("O23").Value = 40
$ws.Range("K24").Value = 32
$ws.Range("L24").Value = 24
$ws.Range("N24").Value = 16
$ws.Range("O24").Value = 42
$ws.Range("K25").Value = 33
$ws.Range("L25").Value = 26
$ws.Range("N25").Value = 19
$ws.Range("O25").Value = 43
$ws.Range("K26").Value = 35
$ws.Range("L26").Value = 29
$ws.Range("N26").Value = 21
$ws.Range("O26").Value = 46
